$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 value
$ws.Range("B2").Value = 3051

# Update A3 and B3 values
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 1207

# Delete rows 4 and 5 (their data is removed entirely per diff)
$ws.Range("A4:B5").Delete()
